# "Generate Report for Handback"
#
# The localization-status workbook gets a handback pass over it:
#   - the Overview/status text "Ready for handoff" becomes
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview!B2:C3 and the Status column on each language sheet)
#   - each language sheet (zh-cn, de-de) grows two new populated
#     columns: F "Latest Target File" and G "Latest Handback File",
#     each a hyperlinked file name, for both data rows (2 and 3)
#   - the "Latest Handback DateTime" column (H) moves on from the
#     zero-date placeholder to the actual handback timestamp

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$srcRepoSha = "6d70edb1aed24ea26c94f1fb6634217c55bca4e7"

# ---- Overview sheet: refresh the status text -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- Per-language sheet handback report -----------------------------
# (sheet name, xlf display name, handoff-repo sha, handback datetime)
$langs = @(
    @{
        Sheet   = "zh-cn"
        Xlf     = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandoffSha = "3b0ab8d59bd64cb2a626a2b6143c0682a7388d7a"
        HandoffOrg = "oltest.zh-cn"
        Handback = "2016-03-19 06:26:23"
    },
    @{
        Sheet   = "de-de"
        Xlf     = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandoffSha = "657f1b63933d736a13bf1584c7eaafc5d4964919"
        HandoffOrg = "oltest.de-de"
        Handback = "2016-03-19 06:26:28"
    }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both data rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $lang.HandoffSha + "/ol-handoff/OpenLocalizationTestOrg/" + $lang.HandoffOrg + "/ci/ht/" + $lang.Xlf
    $mdUrlA = "https://github.com/OpenLocalizationTest/oltest/blob/" + $srcRepoSha + "/e2e/a.md"

    # Row 2 ("a.md"): populate F (Latest Target File) and G (Latest Handback File)
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrlA, "", "", "a.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", $lang.Xlf) | Out-Null

    # Row 3 ("b.md"): same Target/Handback file pairing as row 2
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlA, "", "", "a.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, "", "", $lang.Xlf) | Out-Null

    # Latest Handback DateTime (H) now reflects the real handback time
    $ws.Range("H2").Value = $lang.Handback
    $ws.Range("H3").Value = $lang.Handback
}

Write-Output "Handback report generated."
